$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 101, shifting existing rows 101:219 down to 102:220.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with its data.
$ws.Cells.Item(101, 1).Value = 1
$ws.Cells.Item(101, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(101, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(101, 4).Value = 44664
$ws.Cells.Item(101, 5).Value = 15
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value = 100108006
$ws.Cells.Item(101, 10).Value = "Plátano"
$ws.Cells.Item(101, 11).Value = "Sin especificar"
$ws.Cells.Item(101, 12).Value = "Pintón"
$ws.Cells.Item(101, 13).Value = 120
$ws.Cells.Item(101, 14).Value = 16000
$ws.Cells.Item(101, 15).Value = 17000
$ws.Cells.Item(101, 16).Value = 16500
$ws.Cells.Item(101, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(101, 18).Value = "Ecuador"
$ws.Cells.Item(101, 19).Value = 825
$ws.Cells.Item(101, 20).Value = 20
